$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'30.023.03"
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = "'  +8.44%  "
$ws.Range("E2").Style = $refStyle
$ws.Range("D3").Value = "'1.869.50"
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = "'  +6.22%  "
$ws.Range("E3").Style = $refStyle
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = $refStyle
$ws.Range("E4").Value = "'  +0.35%  "
$ws.Range("E4").Style = $refStyle
$ws.Range("D5").Value = "'249.19"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "'  +2.72%  "
$ws.Range("E5").Style = $refStyle
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "'  +0.37%  "
$ws.Range("E6").Style = $refStyle
$ws.Range("D7").Value = "'0.4969"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "'  +2.69%  "
$ws.Range("E7").Style = $refStyle
$ws.Range("D8").Value = "'45.39"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "'  +9.56%  "
$ws.Range("E8").Style = $refStyle
$ws.Range("D9").Value = "'0.2833"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "'  +7.08%  "
$ws.Range("E9").Style = $refStyle
$ws.Range("D10").Value = "'0.06538"
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = "'  +5.67%  "
$ws.Range("E10").Style = $refStyle
$ws.Range("D11").Value = "'1.870.01"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "'  +6.34%  "
$ws.Range("E11").Style = $refStyle
$ws.Range("D12").Value = "'16.96"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "'  +4.35%  "
$ws.Range("E12").Style = $refStyle
$ws.Range("D13").Value = "'0.07206"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "'  +3.62%  "
$ws.Range("E13").Style = $refStyle
$ws.Range("D14").Value = "'0.6593"
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = "'  +8.03%  "
$ws.Range("E14").Style = $refStyle
$ws.Range("D15").Value = "'84.75"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "'  +8.97%  "
$ws.Range("E15").Style = $refStyle
$ws.Range("D16").Value = "'4.794"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "'  +5.83%  "
$ws.Range("E16").Style = $refStyle
$ws.Range("D17").Value = "'29.998.35"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "'  +8.41%  "
$ws.Range("E17").Style = $refStyle
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "'  +0.43%  "
$ws.Range("E18").Style = $refStyle
$ws.Range("D19").Value = "'12.80"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "'  +10.57%  "
$ws.Range("E19").Style = $refStyle
$ws.Range("D20").Value = "'0.000007485"
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = "'  +5.35%  "
$ws.Range("E20").Style = $refStyle
$ws.Range("D21").Value = "'0.9992"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "'  +0.28%  "
$ws.Range("E21").Style = $refStyle
$ws.Range("D22").Value = "'2.111.91"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "'  +6.92%  "
$ws.Range("E22").Style = $refStyle
$ws.Range("D23").Value = "'4.735"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "'  +5.11%  "
$ws.Range("E23").Style = $refStyle
$ws.Range("D24").Value = "'9.030"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "'  +6.48%  "
$ws.Range("E24").Style = $refStyle
$ws.Range("D25").Value = "'5.492"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "'  +6.84%  "
$ws.Range("E25").Style = $refStyle
$ws.Range("D26").Value = "'144.73"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "'  +2.37%  "
$ws.Range("E26").Style = $refStyle
$ws.Range("D27").Value = "'135.40"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "'  +24.13%  "
$ws.Range("E27").Style = $refStyle
$ws.Range("D28").Value = "'16.73"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "'  +8.63%  "
$ws.Range("E28").Style = $refStyle
$ws.Range("D29").Value = "'1.943"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "'  +3.71%  "
$ws.Range("E29").Style = $refStyle
$ws.Range("D30").Value = "'1.391"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "'  +0.64%  "
$ws.Range("E30").Style = $refStyle
$ws.Range("D31").Value = "'4.233"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "'  +6.16%  "
$ws.Range("E31").Style = $refStyle
$ws.Range("D32").Value = "'0.08620"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "'  +6.69%  "
$ws.Range("E32").Style = $refStyle
$ws.Range("D33").Value = "'3.879"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "'  +4.50%  "
$ws.Range("E33").Style = $refStyle
$ws.Range("D34").Value = "'0.05070"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "'  +8.55%  "
$ws.Range("E34").Style = $refStyle
$ws.Range("D35").Value = "'1.127"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "'  +9.96%  "
$ws.Range("E35").Style = $refStyle
$ws.Range("D36").Value = "'0.6838"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "'  +9.34%  "
$ws.Range("E36").Style = $refStyle
$ws.Range("D37").Value = "'2.690"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "'  +3.06%  "
$ws.Range("E37").Style = $refStyle
$ws.Range("D38").Value = "'2.327"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "'  +13.84%  "
$ws.Range("E38").Style = $refStyle
$ws.Range("D39").Value = "'2.736"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "'  +6.56%  "
$ws.Range("E39").Style = $refStyle
$ws.Range("D40").Value = "'0.9606"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "'  +3.43%  "
$ws.Range("E40").Style = $refStyle
$ws.Range("D41").Value = "'0.01628"
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = "'  +7.96%  "
$ws.Range("E41").Style = $refStyle
$ws.Range("D42").Value = "'6.099"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "'  +6.55%  "
$ws.Range("E42").Style = $refStyle
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "'  +0.55%  "
$ws.Range("E43").Style = $refStyle
$ws.Range("D44").Value = "'103.67"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "'  +4.02%  "
$ws.Range("E44").Style = $refStyle
$ws.Range("D45").Value = "'0.4183"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "'  +7.81%  "
$ws.Range("E45").Style = $refStyle
$ws.Range("D46").Value = "'7.420"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "'  +6.79%  "
$ws.Range("E46").Style = $refStyle
$ws.Range("D47").Value = "'0.1254"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "'  +7.63%  "
$ws.Range("E47").Style = $refStyle
$ws.Range("D48").Value = "'0.05618"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "'  +4.57%  "
$ws.Range("E48").Style = $refStyle
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("B49").Style = $refStyle
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").Style = $refStyle
$ws.Range("D49").Value = "'8.297"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "'  +4.63%  "
$ws.Range("E49").Style = $refStyle
$ws.Range("B50").Value = "'Elrond"
$ws.Range("B50").Style = $refStyle
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C50").Style = $refStyle
$ws.Range("D50").Value = "'32.40"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "'  +7.53%  "
$ws.Range("E50").Style = $refStyle
$ws.Range("D51").Value = "'0.3727"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "'  +9.40%  "
$ws.Range("E51").Style = $refStyle
